# Update evaluation metrics across the three worksheets with the new
# (final) evaluation run's values.

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.4359430604982206
$wsSummary.Range("C2").Value = 0.07871720116618076
$wsSummary.Range("D2").Value = 0.9642857142857143
$wsSummary.Range("E2").Value = 0.1455525606469003
$wsSummary.Range("F2").Value = 0.2967032967032967
$wsSummary.Range("G2").Value = 0.673058485139022
$wsSummary.Range("H2").Value = 0.8039727126805778
$wsSummary.Range("I2").Value = 27
$wsSummary.Range("J2").Value = 316
$wsSummary.Range("K2").Value = 218
$wsSummary.Range("L2").Value = 1

# --- Classification Report sheet ---
$wsReport = $wb.Worksheets.Item("Classification Report")
$wsReport.Range("B2").Value = 0.9954337899543378
$wsReport.Range("C2").Value = 0.4082397003745318
$wsReport.Range("D2").Value = 0.5790172642762285

$wsReport.Range("B3").Value = 0.07871720116618076
$wsReport.Range("C3").Value = 0.9642857142857143
$wsReport.Range("D3").Value = 0.1455525606469003

$wsReport.Range("B4").Value = 0.4359430604982206
$wsReport.Range("C4").Value = 0.4359430604982206
$wsReport.Range("D4").Value = 0.4359430604982206
$wsReport.Range("E4").Value = 0.4359430604982206

$wsReport.Range("B5").Value = 0.5370754955602594
$wsReport.Range("C5").Value = 0.6862627073301231
$wsReport.Range("D5").Value = 0.3622849124615644

$wsReport.Range("B6").Value = 0.9497610773456752
$wsReport.Range("C6").Value = 0.4359430604982206
$wsReport.Range("D6").Value = 0.5574211580455858

# --- Confusion Matrix sheet ---
$wsConf = $wb.Worksheets.Item("Confusion Matrix")
$wsConf.Range("B2").Value = 218
$wsConf.Range("C2").Value = 316
$wsConf.Range("B3").Value = 1
$wsConf.Range("C3").Value = 27
